$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the H1
#    heading "Play Cold Spell Free: Review and Top Bonuses".
# ------------------------------------------------------------------
$heading = $d.Paragraphs.Item(1)
$heading.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs.Item(2)
$metaPara.Style = "Normal"

# Bold "Meta description" run
$boldRange = $metaPara.Range
$boldRange.Collapse(1)
$boldRange.InsertAfter("Meta description")
$boldRange.Bold = 1

# Plain run with the rest of the sentence, inserted right before the
# paragraph mark so it stays inside this paragraph.
$tailRange = $d.Range($metaPara.Range.End - 1, $metaPara.Range.End - 1)
$tailRange.InsertAfter(": Experience the icy realm of Cold Spell with our review. Play for free and land the Fire Wizard Scatter symbol to earn big rewards.")
$tailRange.Bold = 0

# ------------------------------------------------------------------
# 2) Near the end of the document, remove the duplicated bold
#    "Play Cold Spell Free: Review and Top Bonuses" paragraph and turn
#    the following italic paragraph into the image-prompt text.
# ------------------------------------------------------------------
$found = $false
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text.TrimEnd("`r")
    if ($t -eq "Play Cold Spell Free: Review and Top Bonuses") {
        $p.Range.Delete()
        $found = $true
        break
    }
}

if ($found) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        $t = $p.Range.Text.TrimEnd("`r")
        if ($t -eq "Experience the icy realm of Cold Spell with our review. Play for free and land the Fire Wizard Scatter symbol to earn big rewards.") {
            $textRange = $d.Range($p.Range.Start, $p.Range.End - 1)
            $textRange.Text = 'Please create a cartoon-style feature image for the online slot game "Cold Spell". The image should feature a happy Maya warrior with glasses. The Maya warrior should be depicted holding a wand and standing in front of ice-covered mountains with a snowy background to reflect the game''s medieval fantasy theme. The image can include other elements from the game such as playing cards, tiaras, maps, and treasure chests. The image should be bright and colorful to capture the attention of players and entice them to try the game.'
            break
        }
    }
}

Write-Host "Done"
